$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 16838.34481048651
$ws.Range("R2").Value = 151545.1032943786
$ws.Range("S2").Value = 0.19394309853408
$ws.Range("T2").Value = 0.19394309853408

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 15910.63610321822
$ws.Range("R3").Value = 143195.724928964
$ws.Range("S3").Value = 0.1832578023693046
$ws.Range("T3").Value = 0.1832578023693047

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 149.656361
$ws.Range("H4").Value = 448.969083
$ws.Range("I4").Value = 0.5921360794347563
$ws.Range("J4").Value = 0.5921360794347564
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 18660.89938425291
$ws.Range("R4").Value = 167948.0944582762
$ws.Range("S4").Value = 0.2149351785313717
$ws.Range("T4").Value = 0.2149351785313718

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 5933.295240397536
$ws.Range("R5").Value = 53399.65716357782
$ws.Range("S5").Value = 0.0683393573650758
$ws.Range("T5").Value = 0.0683393573650758

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 52.73412466666667
$ws.Range("H6").Value = 158.202374
$ws.Range("I6").Value = 0.208649853730866
$ws.Range("J6").Value = 0.208649853730866
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 5606.400303914094
$ws.Range("R6").Value = 50457.60273522685
$ws.Range("S6").Value = 0.0645742000654571
$ws.Range("T6").Value = 0.06457420006545712

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 52.73412466666667
$ws.Range("H7").Value = 158.202374
$ws.Range("I7").Value = 0.208649853730866
$ws.Range("J7").Value = 0.208649853730866
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 6575.505297241032
$ws.Range("R7").Value = 59179.54767516928
$ws.Range("S7").Value = 0.07573629630033309
$ws.Range("T7").Value = 0.07573629630033311

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.34932566666667
$ws.Range("H8").Value = 151.047977
$ws.Range("I8").Value = 0.1992140668343777
$ws.Range("J8").Value = 0.1992140668343777
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 5664.973415669328
$ws.Range("R8").Value = 50984.76074102395
$ws.Range("S8").Value = 0.06524884183770055
$ws.Range("T8").Value = 0.06524884183770056

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.34932566666667
$ws.Range("H9").Value = 151.047977
$ws.Range("I9").Value = 0.1992140668343777
$ws.Range("J9").Value = 0.1992140668343777
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 5352.86167171176
$ws.Range("R9").Value = 48175.75504540584
$ws.Range("S9").Value = 0.06165395650940461
$ws.Range("T9").Value = 0.06165395650940463

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.34932566666667
$ws.Range("H10").Value = 151.047977
$ws.Range("I10").Value = 0.1992140668343777
$ws.Range("J10").Value = 0.1992140668343777
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 6278.140762293754
$ws.Range("R10").Value = 56503.26686064379
$ws.Range("S10").Value = 0.0723112684872725
$ws.Range("T10").Value = 0.07231126848727251
